$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 5 through 14 (old 2003-2012 data rows), leaving rows 2-4 (2000-2002) to be overwritten below
$ws.Range("A5:E14").EntireRow.Delete()

# Overwrite rows 2-4 with the last three years of data (2010-2012)
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 53243.9321
$ws.Range("C2").Value = 24670.3854
$ws.Range("D2").Value = 28573.5467
$ws.Range("E2").Value = 27694.7704

$ws.Range("A3").Value = "2011年"
$ws.Range("B3").Value = 53685.4444
$ws.Range("C3").Value = 24933.3797
$ws.Range("D3").Value = 28752.0545
$ws.Range("E3").Value = 27355.4198

$ws.Range("A4").Value = "2012年"
$ws.Range("B4").Value = 53857.8772
$ws.Range("C4").Value = 25009.9082
$ws.Range("D4").Value = 28847.0026
$ws.Range("E4").Value = 27032.2501
